$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 195, shifting the old
# rows 195-219 down to become rows 197-221 (preserves all their data).
$ws.Rows.Item(195).Insert()
$ws.Rows.Item(195).Insert()

# Populate the first new row (195) - "Primera" quality entry for 2021-09-10
$ws.Range("A195").Value = 9
$ws.Range("B195").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C195").Value = "Metropolitana"
$ws.Range("D195").Value = 44449
$ws.Range("E195").Value = 13
$ws.Range("F195").Value = 100112012
$ws.Range("G195").Value = "Espinaca"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 240
$ws.Range("K195").Value = 7000
$ws.Range("L195").Value = 8000
$ws.Range("M195").Value = 7500
$ws.Range("N195").Value = "$/cuna 10 kilos"
$ws.Range("O195").Value = "Provincia de Chacabuco"
$ws.Range("P195").Value = 750
$ws.Range("Q195").Value = 10
$ws.Range("R195").Value = "Hortaliza"

# Populate the second new row (196) - "Segunda" quality entry for 2021-09-10
$ws.Range("A196").Value = 9
$ws.Range("B196").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C196").Value = "Metropolitana"
$ws.Range("D196").Value = 44449
$ws.Range("E196").Value = 13
$ws.Range("F196").Value = 100112012
$ws.Range("G196").Value = "Espinaca"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Segunda"
$ws.Range("J196").Value = 115
$ws.Range("K196").Value = 5000
$ws.Range("L196").Value = 6000
$ws.Range("M196").Value = 5496
$ws.Range("N196").Value = "$/cuna 10 kilos"
$ws.Range("O196").Value = "Provincia de Chacabuco"
$ws.Range("P196").Value = 550
$ws.Range("Q196").Value = 10
$ws.Range("R196").Value = "Hortaliza"
